$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.299.53"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.690.28"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'217.68"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'0.5395"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.2731"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'0.06444"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'21.66"
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.07667"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Value = "1.718.26"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "'4.538"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "'0.5785"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'0.000008378"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "'66.89"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "26.366.75"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'4.904"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'10.86"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'190.50"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'6.254"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'149.06"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").Value = "'0.1285"
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("D26").Value = "'7.857"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'15.88"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'0.06284"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").Value = "'1.372"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'3.600"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'3.579"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "'1.676"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'0.6165"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").Value = "'2.416"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "'0.01651"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "1.109.95"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'6.113"
$ws.Range("E40").Value = "  -5.29%  "
$ws.Range("D41").Value = "'0.8821"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'1.013"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'101.17"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "1.843.30"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'57.63"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'8.154"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D51").Value = "'6.046"
$ws.Range("E51").Value = "  -0.40%  "
